# KRA master database update
# - Adds three new tracking columns (Merged_From_Count, Merge_Sources, Best_Score)
#   to the KRA_Database sheet and backfills them for the existing rows.
# - Appends eight freshly-merged records (rows 4-11).
# - Refreshes the Database_Summary statistics sheet to reflect the new totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KRA_Database")
$summary = $wb.Worksheets.Item("Database_Summary")

# ---------------------------------------------------------------------------
# 1. New header columns L:N on row 1 - copy formatting from the existing
#    bold/bordered header (column K) so the new headers match the others.
# ---------------------------------------------------------------------------
$ws.Range("L1").Value = "Merged_From_Count"
$ws.Range("M1").Value = "Merge_Sources"
$ws.Range("N1").Value = "Best_Score"

$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Backfill the new columns for the existing two records.
# ---------------------------------------------------------------------------
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = "Unknown, Unknown"
$ws.Range("N2").Value = 100

# Row 3 (the post-migration test row) did not take part in any merge,
# so L3/M3/N3 stay blank (no value to assign).

# Row 3's "year" value was stored as text - normalize it to a real number.
$ws.Range("F3").Value = 2024

# ---------------------------------------------------------------------------
# 3. Append the eight newly extracted/merged records (rows 4-11).
#    Columns D (preAmount) and F (year) hold digit-only strings that must
#    stay text, so force a text number format on those cells before writing.
# ---------------------------------------------------------------------------
$ws.Range("D4:D11").NumberFormat = "@"
$ws.Range("F4:F11").NumberFormat = "@"

$rows = @(
    @{ Row=4;  A="04TH September, 2025"; B="A018905312S"; C="Daisy Jepkosgei Biwott";   D="2025."; F="2024"; G="Franciscar Nyangweta"; H="KITALE";   K=3;  L=2; M="Unknown, Unknown"; N=94.8 },
    @{ Row=5;  A="4th September, 2025";  B="A009775891W"; C="Ezekiel Kipserem Korir";   D="2025."; F="2024"; G="Franciscar Nyangweta"; H="KITALE";   K=4;  L=2; M="Unknown, Unknown"; N=94.8 },
    @{ Row=6;  A="04th September, 2025"; B="A004578892U"; C="JESSY KAGONDU WAMBUGU";    D="2025."; F="2024"; G="Franciscar Nyangweta"; H="KITALE";   K=6;  L=2; M="Unknown, Unknown"; N=94.8 },
    @{ Row=7;  A="04th September, 2025"; B="A008596925K"; C="KELVIN KIPKEMBOI MUTAI";   D="2025."; F="2024"; G="Franciscar Nyangweta"; H="KITALE";   K=7;  L=2; M="Unknown, Unknown"; N=94.8 },
    @{ Row=8;  A="10th September, 2025"; B="A007388222W"; C="MICHAEL MWANGI MUCHUNGI";  D="2025."; F="2024"; G="Franciscar Nyangweta"; H="KITALE";   K=8;  L=2; M="Unknown, Unknown"; N=94.8 },
    @{ Row=9;  A="04th September, 2025"; B="A012209532N"; C="Paul Chotomolo Mirikwa";   D="2025."; F="2024"; G="Franciscar Nyangweta"; H="NAITIRI";  K=9;  L=2; M="Unknown, Unknown"; N=94.8 },
    @{ Row=10; A="29TH AUGUST, 2025";    B="A001126762Z"; C="Peter Kimutai Telengech";  D="2025."; F="2024"; G="Franciscar Nyangweta"; H="ELDORET";  K=10; L=2; M="Unknown, Unknown"; N=94.8 },
    @{ Row=11; A="04th September, 2025"; B="A005615142S"; C="THOMAS JUMA SIKUKU";       D="2025";  F="2024"; G="Franciscar Nyangweta"; H="KITALE";   K=11; L=2; M="Unknown, Unknown"; N=94.8 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.A
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    # E (finalAmount) has no data for these records - leave blank.
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
    $ws.Range("H$n").Value = $r.H
    $ws.Range("I$n").Value = "2025-09-22 11:25:52"
    $ws.Range("J$n").Value = "multi_format_extractor"
    $ws.Range("K$n").Value = $r.K
    $ws.Range("L$n").Value = $r.L
    $ws.Range("M$n").Value = $r.M
    $ws.Range("N$n").Value = $r.N
}

# ---------------------------------------------------------------------------
# 4. Refresh the Database_Summary sheet with the new aggregate statistics.
# ---------------------------------------------------------------------------
$summary.Range("B2").Value = 10
$summary.Range("B3").Value = "2025-09-22 11:25:52"
$summary.Range("B4").Value = 9
$summary.Range("B5").Value = 1
$summary.Range("B6").Value = "04TH September, 2025"
$summary.Range("B7").Value = "4th September, 2025"
$summary.Range("B8").Value = 10
$summary.Range("B9").Value = 4

Write-Host "KRA_Database and Database_Summary updated."
